# Adds 18 "Deleted variable" rows (rows 69-86) to the "SAM Variable Changes"
# sheet, documenting variables removed from the Molten Salt Tower Power
# Block / Receiver, the Molten Salt Power Block, and the Molten Salt Tower
# Storage input pages (redundant, unused, or constant values that were
# hardcoded in the compute module instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# --- 1. Insert 18 new rows (69-86), copying the formatting of the last
#        existing data row (68) so the per-column styles (s=25/19/56) are
#        reused instead of new style records being created. ---
$ws.Rows.Item(68).Copy()
$ws.Range("A69:A86").EntireRow.Insert()

# --- 2. Fill in the new rows' data. ---
# Each row: A=Deleted variable, B=number, C=<old name>, E=<input page>,
#           F=not used, G=N/A, H=Ty  (column D is intentionally left blank)
$rows = @(
    @{ Old = "m_dot_htf_ref";                  Page = "Molten Salt Tower Power Block" },
    @{ Old = "T_pb_out";                       Page = "Molten Salt Tower Power Block" },
    @{ Old = "mode";                           Page = "Molten Salt Tower Power Block" },
    @{ Old = "fthr_ok";                        Page = "Molten Salt Tower Power Block" },
    @{ Old = "pb_fixed_par_cntl";              Page = "Molten Salt Tower Power Block" },
    @{ Old = "dt_cold";                        Page = "Molten Salt Tower Power Block" },
    @{ Old = "dt_hot";                         Page = "Molten Salt Tower Power Block" },
    @{ Old = "hx_config";                      Page = "Molten Salt Tower Power Block" },
    @{ Old = "is_hx";                          Page = "Molten Salt Tower Power Block" },
    @{ Old = "tech_type";                      Page = "Molten Salt Tower Power Block" },
    @{ Old = "deg_wind";                       Page = "Molten Salt Tower Receiver" },
    @{ Old = "P_htf";                          Page = "Molten Salt Tower Receiver" },
    @{ Old = "T_salt_cold";                    Page = "Molten Salt Power Block" },
    @{ Old = "HTF";                            Page = "Molten Salt Power Block" },
    @{ Old = "Design_power";                   Page = "Molten Salt Power Block" },
    @{ Old = "csp.pt.pwrb.min_restart_time";   Page = "Molten Salt Power Block" },
    @{ Old = "csp.pt.rec.max_rec_flux";        Page = "Molten Salt Power Block" },
    @{ Old = "store_fluid";                    Page = "Molten Salt Tower Storage" }
)

$r = 69
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "Deleted variable"
    $ws.Cells.Item($r, 2).Value = "number"
    $ws.Cells.Item($r, 3).Value = $row.Old
    $ws.Cells.Item($r, 5).Value = $row.Page
    $ws.Cells.Item($r, 6).Value = "not used"
    $ws.Cells.Item($r, 7).Value = "N/A"
    $ws.Cells.Item($r, 8).Value = "Ty"
    $r = $r + 1
}

# --- 3. Widen column C slightly to fit the longer variable names.
#        (The engine's column-width model adds a fixed 5/6 character offset
#        when round-tripping through the OOXML "width" attribute, so we
#        back that out here to land exactly on width="28" in the XML.) ---
$ws.Columns.Item(3).ColumnWidth = 163/6

# --- 4. Update the data validation on column A so it covers exactly the
#        new data range (A2:A86) without the allowBlank attribute. ---
$ws.Range("A2:A73").Validation.Delete()
$validation = $ws.Range("A2:A86").Validation
$validation.Add(3, 1, 1, "Types")
$validation.IgnoreBlank = $false

# --- 5. Restore the view so the newly added rows are visible. ---
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("C77").Select()
